$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: Actual Results (I) and Test Results (J) for TC001
$ws.Cells.Item(7, 9).Value = "same as expected"
$ws.Cells.Item(7, 10).Value = "pass"

# New column L header (row 6)
$ws.Cells.Item(6, 12).Value = "Screenshots"

# Row 8: Actual Results (I) and Test Results (J) for step 2
$ws.Cells.Item(8, 9).Value = "same as expected"
$ws.Cells.Item(8, 10).Value = "pass"

# Row 8: hyperlink to screenshot in column L
$ws.Hyperlinks.Add($ws.Cells.Item(8, 12), "D:\010LiveTech\Screenshots\TC001\Screenshot 2024-06-24 100610.png") | Out-Null

# Row 9: Actual Results (I) and Test Results (J)
$ws.Cells.Item(9, 9).Value = "same as expected"
$ws.Cells.Item(9, 10).Value = "pass"

$ws.Columns.Item(9).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(12).EntireColumn.AutoFit() | Out-Null

$ws.Range("L8").Select()
